$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "Mean mass flow rate" simulation output values (simplify data added to simOut)
$ws.Range("B2").Value = 0.03403972726862195
$ws.Range("B3").Value = 0.25977738384266158
$ws.Range("B4").Value = 0.018982081937515962
$ws.Range("B5").Value = 0.031488317742870427
$ws.Range("B6").Value = 0.3516024905798032
$ws.Range("B7").Value = 0.020017272644329322
$ws.Range("B8").Value = 0.032667304614166892
$ws.Range("B9").Value = 0.10700075746217552

# Column widths were re-derived for the new data (best achievable via ColumnWidth)
$ws.Columns("A").ColumnWidth = 23.166666666666668
$ws.Columns("B").ColumnWidth = 17
$ws.Columns("C").ColumnWidth = 3.8333333333333335
